$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('B2').Value = 'Bitcoin'
$ws.Range('C2').Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range('D2').Value = '28.014.76'
$ws.Range('E2').Value = '  -0.21%  '

$ws.Range('B3').Value = 'Ethereum'
$ws.Range('C3').Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range('D3').Value = '1.868.94'
$ws.Range('E3').Value = '  -1.25%  '

$ws.Range('B4').Value = 'TetherUSD'
$ws.Range('C4').Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '1.004'
$c.Style = 'Normal'
$ws.Range('E4').Value = '  +0.30%  '

$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '312.50'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -0.54%  '

$ws.Range('B6').Value = 'USDC'
$ws.Range('C6').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '1.003'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +0.21%  '

$ws.Range('B7').Value = 'XRP'
$ws.Range('C7').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.5100'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  +1.33%  '

$ws.Range('B8').Value = 'Cardano'
$ws.Range('C8').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.3798'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  -2.57%  '

$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.08281'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  -10.47%  '

$ws.Range('B10').Value = 'Polygon'
$ws.Range('C10').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '1.109'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -1.85%  '

$ws.Range('B11').Value = 'OKB'
$ws.Range('C11').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '41.34'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -1.10%  '

$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '6.207'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  -2.76%  '

$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.866.68'
$ws.Range('E13').Value = '  -1.61%  '

$ws.Range('B14').Value = 'Solana'
$ws.Range('C14').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '20.46'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -1.82%  '

$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '7.178'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  -1.73%  '

$ws.Range('B16').Value = 'BinanceUSD'
$ws.Range('C16').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '1.004'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  +0.32%  '

$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '0.00001095'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  -1.30%  '

$ws.Range('B18').Value = 'Litecoin'
$ws.Range('C18').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '90.66'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  -1.38%  '

$ws.Range('B19').Value = 'TRON'
$ws.Range('C19').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '0.06619'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -0.17%  '

$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '17.83'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -0.31%  '

$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '1.002'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +0.17%  '

$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '5.997'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -3.79%  '

$ws.Range('B23').Value = 'WrappedBTC'
$ws.Range('C23').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D23').Value = '28.053.81'
$ws.Range('E23').Value = '  -0.31%  '

$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '11.09'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -2.55%  '

$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '2.259'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -2.79%  '

$ws.Range('B26').Value = 'LidoDAOToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '2.570'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +0.66%  '

$ws.Range('B27').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C27').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D27').Value = '2.085.84'
$ws.Range('E27').Value = '  -1.34%  '

$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '157.23'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  -0.86%  '

$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '20.50'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  -1.53%  '

$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '125.38'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -1.26%  '

$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '0.1057'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +0.14%  '

$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '1.041'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -3.17%  '

$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '5.589'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -0.28%  '

$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '3.595'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -0.17%  '

$ws.Range('B35').Value = 'FraxShare'
$ws.Range('C35').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '9.697'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +2.18%  '

$ws.Range('B36').Value = 'VeChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '0.02433'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +0.75%  '

$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.06536'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -0.87%  '

$ws.Range('B38').Value = 'Algorand'
$ws.Range('C38').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.2162'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -1.72%  '

$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '1.205'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -1.01%  '

$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.6425'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -0.51%  '

$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '1.234'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -8.30%  '

$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '11.27'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -2.96%  '

$ws.Range('B43').Value = 'InternetComputer(DFINITY)'
$ws.Range('C43').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '4.868'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -1.83%  '

$ws.Range('B44').Value = 'Decentraland'
$ws.Range('C44').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '0.6107'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +0.67%  '

$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '13.06'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -2.78%  '

$ws.Range('B46').Value = 'WEMIXTOKEN'
$ws.Range('C46').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '1.288'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -0.88%  '

$ws.Range('B47').Value = 'PancakeSwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '3.658'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -0.93%  '

$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '1.997'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -0.38%  '

$ws.Range('B49').Value = 'EOS'
$ws.Range('C49').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '1.211'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  +1.42%  '

$ws.Range('B50').Value = 'Quant'
$ws.Range('C50').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '121.05'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -0.58%  '

$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '79.90'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +0.84%  '
